# we go. Hard 추가
# Rename the sheet to match the note-map song title, toggle the "Hard"
# chart's split markers (columns D/E/F/G) to their new on/off states, and
# restore the editor's last scroll/zoom/selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename Sheet1 -> WeGo_Fromis9 ---
$ws.Name = "WeGo_Fromis9"

# --- Toggle split markers ---
# Reference cells used purely as format donors for Copy/PasteSpecial so the
# existing shared cell style ("highlighted") is reused instead of a new one
# being allocated in styles.xml.
$styledSrc = $ws.Range("A8")   # already uses the highlighted style
$plainSrc = $ws.Range("D8")    # uses the default (no) style

# Cells turned ON (now marked present for this split): apply highlighted style, value 1
$styledSrc.Copy()
$ws.Range("F278").PasteSpecial(-4122)
$ws.Range("F278").Value = 1
$styledSrc.Copy()
$ws.Range("G304").PasteSpecial(-4122)
$ws.Range("G304").Value = 1
$styledSrc.Copy()
$ws.Range("E310").PasteSpecial(-4122)
$ws.Range("E310").Value = 1
$styledSrc.Copy()
$ws.Range("F336").PasteSpecial(-4122)
$ws.Range("F336").Value = 1
$styledSrc.Copy()
$ws.Range("G360").PasteSpecial(-4122)
$ws.Range("G360").Value = 1
$styledSrc.Copy()
$ws.Range("F362").PasteSpecial(-4122)
$ws.Range("F362").Value = 1
$styledSrc.Copy()
$ws.Range("G368").PasteSpecial(-4122)
$ws.Range("G368").Value = 1
$styledSrc.Copy()
$ws.Range("F370").PasteSpecial(-4122)
$ws.Range("F370").Value = 1
$styledSrc.Copy()
$ws.Range("G376").PasteSpecial(-4122)
$ws.Range("G376").Value = 1
$styledSrc.Copy()
$ws.Range("F378").PasteSpecial(-4122)
$ws.Range("F378").Value = 1
$styledSrc.Copy()
$ws.Range("F440").PasteSpecial(-4122)
$ws.Range("F440").Value = 1
$styledSrc.Copy()
$ws.Range("E448").PasteSpecial(-4122)
$ws.Range("E448").Value = 1
$styledSrc.Copy()
$ws.Range("F472").PasteSpecial(-4122)
$ws.Range("F472").Value = 1
$styledSrc.Copy()
$ws.Range("E480").PasteSpecial(-4122)
$ws.Range("E480").Value = 1
$styledSrc.Copy()
$ws.Range("F488").PasteSpecial(-4122)
$ws.Range("F488").Value = 1
$styledSrc.Copy()
$ws.Range("F492").PasteSpecial(-4122)
$ws.Range("F492").Value = 1
$styledSrc.Copy()
$ws.Range("E496").PasteSpecial(-4122)
$ws.Range("E496").Value = 1
$styledSrc.Copy()
$ws.Range("E500").PasteSpecial(-4122)
$ws.Range("E500").Value = 1
$styledSrc.Copy()
$ws.Range("G518").PasteSpecial(-4122)
$ws.Range("G518").Value = 1
$styledSrc.Copy()
$ws.Range("F520").PasteSpecial(-4122)
$ws.Range("F520").Value = 1
$styledSrc.Copy()
$ws.Range("D636").PasteSpecial(-4122)
$ws.Range("D636").Value = 1
$styledSrc.Copy()
$ws.Range("F640").PasteSpecial(-4122)
$ws.Range("F640").Value = 1
$styledSrc.Copy()
$ws.Range("E644").PasteSpecial(-4122)
$ws.Range("E644").Value = 1
$styledSrc.Copy()
$ws.Range("D652").PasteSpecial(-4122)
$ws.Range("D652").Value = 1
$styledSrc.Copy()
$ws.Range("E654").PasteSpecial(-4122)
$ws.Range("E654").Value = 1
$styledSrc.Copy()
$ws.Range("D656").PasteSpecial(-4122)
$ws.Range("D656").Value = 1
$styledSrc.Copy()
$ws.Range("D804").PasteSpecial(-4122)
$ws.Range("D804").Value = 1
$styledSrc.Copy()
$ws.Range("E928").PasteSpecial(-4122)
$ws.Range("E928").Value = 1
$styledSrc.Copy()
$ws.Range("E1184").PasteSpecial(-4122)
$ws.Range("E1184").Value = 1

# Cells turned OFF (no longer marked for this split): clear style, value 0
$plainSrc.Copy()
$ws.Range("D236").PasteSpecial(-4122)
$ws.Range("D236").Value = 0
$plainSrc.Copy()
$ws.Range("E238").PasteSpecial(-4122)
$ws.Range("E238").Value = 0
$plainSrc.Copy()
$ws.Range("F360").PasteSpecial(-4122)
$ws.Range("F360").Value = 0
$plainSrc.Copy()
$ws.Range("G366").PasteSpecial(-4122)
$ws.Range("G366").Value = 0
$plainSrc.Copy()
$ws.Range("F368").PasteSpecial(-4122)
$ws.Range("F368").Value = 0
$plainSrc.Copy()
$ws.Range("G374").PasteSpecial(-4122)
$ws.Range("G374").Value = 0
$plainSrc.Copy()
$ws.Range("F376").PasteSpecial(-4122)
$ws.Range("F376").Value = 0
$plainSrc.Copy()
$ws.Range("G382").PasteSpecial(-4122)
$ws.Range("G382").Value = 0
$plainSrc.Copy()
$ws.Range("G492").PasteSpecial(-4122)
$ws.Range("G492").Value = 0
$plainSrc.Copy()
$ws.Range("F500").PasteSpecial(-4122)
$ws.Range("F500").Value = 0
$plainSrc.Copy()
$ws.Range("F636").PasteSpecial(-4122)
$ws.Range("F636").Value = 0
$plainSrc.Copy()
$ws.Range("G640").PasteSpecial(-4122)
$ws.Range("G640").Value = 0
$plainSrc.Copy()
$ws.Range("G644").PasteSpecial(-4122)
$ws.Range("G644").Value = 0
$plainSrc.Copy()
$ws.Range("F652").PasteSpecial(-4122)
$ws.Range("F652").Value = 0
$plainSrc.Copy()
$ws.Range("D802").PasteSpecial(-4122)
$ws.Range("D802").Value = 0
$plainSrc.Copy()
$ws.Range("E926").PasteSpecial(-4122)
$ws.Range("E926").Value = 0
$plainSrc.Copy()
$ws.Range("E1182").PasteSpecial(-4122)
$ws.Range("E1182").Value = 0

# --- Restore view state: scroll position, zoom, and active selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 55
$ws.Range("S1227").Select()

